$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 199, shifting existing rows
# 199-283 down to 200-284.
$ws.Rows(199).Insert()

# Populate the newly inserted row 199 with the new record's data.
$ws.Cells.Item(199, 1).Value = 8
$ws.Cells.Item(199, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(199, 3).Value = "Coquimbo"
$ws.Cells.Item(199, 4).Value = 44636
$ws.Cells.Item(199, 5).Value = 4
$ws.Cells.Item(199, 6).Value = 100112032
$ws.Cells.Item(199, 7).Value = "Zapallo italiano"
$ws.Cells.Item(199, 8).Value = "Sin especificar"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 560
$ws.Cells.Item(199, 11).Value = 11000
$ws.Cells.Item(199, 12).Value = 12000
$ws.Cells.Item(199, 13).Value = 11500
$ws.Cells.Item(199, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(199, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(199, 16).Value = 192
$ws.Cells.Item(199, 17).Value = 60
$ws.Cells.Item(199, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(199, 4).NumberFormat = $ws.Cells.Item(200, 4).NumberFormat
